# Literary criticism - EL2410 addition to sem 3
# Two students' grades were recalculated (one in sem2, one in sem3) which
# changed their CGPA and moved them up in the ranked list. This reorders
# the affected block of rows while keeping everything else the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1: rows 59-65 (ranks 58-64) ---
# Student 230218 (GUNATHUNGA U.A.) had sem2 recalculated (3.775 -> 3.8179),
# raising CGPA from 3.8129 to 3.8323 and jumping from rank 64 to rank 58.
# Everyone previously ranked 58-63 shifts down by one place.

$ws.Range("B59").Value = 230218
$ws.Range("C59").Value = "GUNATHUNGA U.A."
$ws.Range("D59").Value = 3.9357
$ws.Range("E59").Value = 3.8179
$ws.Range("F59").Value = 3.78
$ws.Range("G59").Value = 3.8323

$ws.Range("B60").Value = 230500
$ws.Range("C60").Value = "PRISHMIKA H.W.N."
$ws.Range("D60").Value = 3.9571
$ws.Range("E60").Value = 3.9045
$ws.Range("F60").Value = 3.6783
$ws.Range("G60").Value = 3.8288

$ws.Range("B61").Value = 230629
$ws.Range("C61").Value = "TENNAKOON U.G.R.B."
$ws.Range("D61").Value = 3.9571
$ws.Range("E61").Value = 3.9045
$ws.Range("F61").Value = 3.6739
$ws.Range("G61").Value = 3.8271

$ws.Range("B62").Value = 230470
$ws.Range("C62").Value = "PEIRIS T.S.R."
$ws.Range("D62").Value = 4
$ws.Range("E62").Value = 4
$ws.Range("F62").Value = 3.5261
$ws.Range("G62").Value = 3.8242

$ws.Range("B63").Value = 230180
$ws.Range("C63").Value = "FERNANDO H.M.D."
$ws.Range("D63").Value = 3.9357
$ws.Range("E63").Value = 3.832
$ws.Range("F63").Value = 3.7391
$ws.Range("G63").Value = 3.821

$ws.Range("B64").Value = 230353
$ws.Range("C64").Value = "KUMARA P.K.M.P."
$ws.Range("D64").Value = 3.9
$ws.Range("E64").Value = 3.964
$ws.Range("F64").Value = 3.613
$ws.Range("G64").Value = 3.8194

$ws.Range("B65").Value = 230502
$ws.Range("C65").Value = "PRIYADARSHANA S.A.D."
$ws.Range("D65").Value = 4
$ws.Range("E65").Value = 3.9455
$ws.Range("F65").Value = 3.587
$ws.Range("G65").Value = 3.8186

# --- Block 2: rows 98-100 (ranks 97-99) ---
# Student 230013 (ABEYWARNA D.H.) had sem3 recalculated (3.5913 -> 3.6385,
# the EL2410 addition), raising CGPA from 3.6532 to 3.6692 and jumping
# from rank 99 to rank 97. Everyone previously ranked 97-98 shifts down.

$ws.Range("B98").Value = 230013
$ws.Range("C98").Value = "ABEYWARNA D.H."
$ws.Range("D98").Value = 3.85
$ws.Range("E98").Value = 3.6
$ws.Range("F98").Value = 3.6385
$ws.Range("G98").Value = 3.6692

$ws.Range("B99").Value = 230229
$ws.Range("C99").Value = "HANSINDU M.M.A.D."
$ws.Range("D99").Value = 3.85
$ws.Range("E99").Value = 3.784
$ws.Range("F99").Value = 3.4
$ws.Range("G99").Value = 3.665

$ws.Range("B100").Value = 230650
$ws.Range("C100").Value = "UBEYSEKARA V.T.T."
$ws.Range("D100").Value = 4
$ws.Range("E100").Value = 3.7364
$ws.Range("F100").Value = 3.3739
$ws.Range("G100").Value = 3.6576
